$d = $word.ActiveDocument

# Replace version number ".0" -> ".1" for "Version: 4.2.0" -> "Version: 4.2.1"
$d.Content.Find.Execute("Version: 4.2.0", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Version: 4.2.1", 2)

# Replace the published date "January 2024" -> "March 2024"
$d.Content.Find.Execute("Published date: January 2024", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Published date: March 2024", 2)
